$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.63
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 4.75
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 19
$ws.Range("AC3").Value = 7.5
$ws.Range("AI3").Value = 19
$ws.Range("AK3").Value = 41
$ws.Range("AO3").Value = 11
$ws.Range("AX3").Value = 23
$ws.Range("G5").Value = 2.8
$ws.Range("I5").Value = 2.45
$ws.Range("J5").Value = 3.5
$ws.Range("AL5").Value = 19
$ws.Range("AO5").Value = 17
$ws.Range("AP5").Value = 26
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 6.5
$ws.Range("K6").Value = 2.42
$ws.Range("P6").Value = 4.1
$ws.Range("X6").Value = 6
$ws.Range("AC6").Value = 13.5
$ws.Range("AL6").Value = 50
$ws.Range("AN6").Value = 3.3
$ws.Range("AT6").Value = 3.25
$ws.Range("G7").Value = 2.37
$ws.Range("H7").Value = 3.15
$ws.Range("I7").Value = 2.85
$ws.Range("J7").Value = 2.87
$ws.Range("L7").Value = 3.45
$ws.Range("V7").Value = 2.18
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 14
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 27
$ws.Range("AA7").Value = 17.5
$ws.Range("AB7").Value = 22
$ws.Range("AC7").Value = 11.25
$ws.Range("AD7").Value = 6.3
$ws.Range("AF7").Value = 45
$ws.Range("AH7").Value = 10
$ws.Range("AI7").Value = 16
$ws.Range("AJ7").Value = 10
$ws.Range("AK7").Value = 37
$ws.Range("AL7").Value = 23
$ws.Range("AM7").Value = 27
$ws.Range("AN7").Value = 4.4
$ws.Range("AO7").Value = 12
$ws.Range("AP7").Value = 17
$ws.Range("AQ7").Value = 45
$ws.Range("AR7").Value = 65
$ws.Range("AT7").Value = 2.82
$ws.Range("AV7").Value = 50
$ws.Range("AW7").Value = 4.9
$ws.Range("AX7").Value = 16
$ws.Range("AY7").Value = 21
$ws.Range("AZ7").Value = 75
$ws.Range("BA7").Value = 100
$ws.Range("BB7").Value = 250
